$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns keep text formatting so values like
# "236.18" or "1.00" are not auto-converted into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '96.491.29'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '3.706.17'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '236.18'
$ws.Range("E5").Value = '  -3.28%  '
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("D7").Value = '650.11'
$ws.Range("E7").Value = '  -3.32%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = '1.06'
$ws.Range("E10").Value = '  -6.34%  '
$ws.Range("D11").Value = '3.702.77'
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").Value = '44.27'
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '0.0000305'
$ws.Range("E13").Value = '  +15.68%  '
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").Value = '6.71'
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("D16").Value = '4.392.68'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '96.345.45'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").Value = '8.82'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '3.701.94'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '13.04'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").Value = '18.61'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  -8.61%  '
$ws.Range("D23").Value = '519.90'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("D24").Value = '3.39'
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("D25").Value = '0.0000210'
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("D26").Value = '6.92'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = '101.24'
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("D28").Value = '13.26'
$ws.Range("E28").Value = '  +1.84%  '
$ws.Range("D29").Value = '0.176'
$ws.Range("E29").Value = '  +4.00%  '
$ws.Range("D30").Value = '3.00'
$ws.Range("E30").Value = '  -2.96%  '
$ws.Range("D31").Value = '12.09'
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").Value = '1.87'
$ws.Range("E33").Value = '  +6.98%  '
$ws.Range("D34").Value = '0.186'
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").Value = '32.24'
$ws.Range("E36").Value = '  -3.14%  '
$ws.Range("D37").Value = '646.73'
$ws.Range("E37").Value = '  +4.82%  '
$ws.Range("D38").Value = '0.587'
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").Value = '8.80'
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").Value = '6.84'
$ws.Range("E41").Value = '  +11.53%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '41.15'
$ws.Range("E42").Value = '  -2.93%  '
$ws.Range("B43").Value = 'ImmutableX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D43").Value = '2.05'
$ws.Range("E43").Value = '  +4.66%  '
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("D45").Value = '0.964'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").Value = '0.0452'
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").Value = '0.432'
$ws.Range("E47").Value = '  +1.32%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '2.28'
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").Value = '23.58'
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").Value = '8.48'
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("E51").Value = '  +1.96%  '

